# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the formatting of an existing header cell (bold, centered, bordered)
# onto the three new header cells so they reuse the same style (s="1")
# instead of minting a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows ---------------------------------------------------------
# Every player row gets the team's season record repeated: 95 wins,
# 67 losses, 0 ties.
$wins = 95
$losses = 67
$ties = 0

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}

Write-Output "season record columns added"
